{"js": "// Apply the \"Added many more features\" edits to the Kronos Unleashed\n// review document. Each change is a like-for-like text replacement, so we\n// locate the original text with Body.search() (exact / case-sensitive\n// match) and rewrite the matched range's text in place. This keeps each\n// run's existing formatting (bold/italic/paragraph style) untouched.\n\nconst body = context.document.body;\n\nconst replacements = [\n  // Title heading (also reused verbatim later in the document as a bold run).\n  [\n    \"Play Kronos Unleashed for Free - Review of WMS Slot Game\",\n    \"Play Kronos Unleashed Free - Review of WMS Slot Game\",\n  ],\n  // \"What we like\" bullet list.\n  [\n    \"Lightning Respins feature for increased pay lines\",\n    \"Regular structure with 60 pay lines and 5 reels\",\n  ],\n  [\n    \"Special symbols represent elements of Greek mythology\",\n    \"Lightning Respins feature with up to 100 pay lines\",\n  ],\n  [\n    \"Stunning visuals and graphics\",\n    \"Chance to win free spins with Scatter symbol\",\n  ],\n  [\n    \"Chance to play for free with Scatter symbol\",\n    \"Stunning visuals and graphics with Greek mythology theme\",\n  ],\n  // \"What we don't like\" bullet list.\n  [\n    \"Limited variety of special features\",\n    \"Limited variety in special symbols\",\n  ],\n  [\n    \"Limited appeal to non-Greek mythology fans\",\n    \"May not appeal to players not interested in Greek mythology\",\n  ],\n  // Closing italic summary paragraph.\n  [\n    \"Immerse yourself in the world of Greek mythology with Kronos Unleashed, an enjoyable and visually stunning slot game by WMS. Play for free and explore special features including Lightning Respins and Scatter symbols.\",\n    \"Read our review of Kronos Unleashed and play this WMS slot game for free with Greek mythology theme and special features.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Kronos Unleashed\n# review document. Each change is a like-for-like text replacement, so we\n# drive Word's Find/Replace (wdReplaceAll) across the whole document body\n# for each old/new text pair. Find/Replace rewrites the text of the\n# matched run(s) in place, so existing run/paragraph formatting (bold,\n# italic, styles) is preserved untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Kronos Unleashed for Free - Review of WMS Slot Game\", \"Play Kronos Unleashed Free - Review of WMS Slot Game\"),\n    @(\"Lightning Respins feature for increased pay lines\", \"Regular structure with 60 pay lines and 5 reels\"),\n    @(\"Special symbols represent elements of Greek mythology\", \"Lightning Respins feature with up to 100 pay lines\"),\n    @(\"Stunning visuals and graphics\", \"Chance to win free spins with Scatter symbol\"),\n    @(\"Chance to play for free with Scatter symbol\", \"Stunning visuals and graphics with Greek mythology theme\"),\n    @(\"Limited variety of special features\", \"Limited variety in special symbols\"),\n    @(\"Limited appeal to non-Greek mythology fans\", \"May not appeal to players not interested in Greek mythology\"),\n    @(\"Immerse yourself in the world of Greek mythology with Kronos Unleashed, an enjoyable and visually stunning slot game by WMS. Play for free and explore special features including Lightning Respins and Scatter symbols.\", \"Read our review of Kronos Unleashed and play this WMS slot game for free with Greek mythology theme and special features.\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
